# Add the 2021 row (row 12) to Sheet1, mirroring the existing 2011-2020 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last data row (A11, which carries the bold /
# bordered / centered "year" style) down onto the new A12 label cell before
# writing its value, so the new row matches the styling of the previous ones.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("A12").Value = "2021年"

$ws.Range("B12").Value = 185.7
$ws.Range("E12").Value = 84715.7
$ws.Range("F12").Value = 7707.7
$ws.Range("G12").Value = 1876.2
$ws.Range("H12").Value = 58413.5
$ws.Range("I12").Value = 71443.10000000001
$ws.Range("J12").Value = 17646.9
$ws.Range("L12").Value = 84715.7
$ws.Range("N12").Value = 82839.5
$ws.Range("S12").Value = 8469.5
$ws.Range("V12").Value = 5564.9

# Columns C, D, K, M, O, P, Q, R, T, U have no reported value for 2021,
# matching the blanks already present in the preceding rows.
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("O12").Value = ""
$ws.Range("P12").Value = ""
$ws.Range("Q12").Value = ""
$ws.Range("R12").Value = ""
$ws.Range("T12").Value = ""
$ws.Range("U12").Value = ""
